# MDS-3438: updated docs footer
#
# The primary footer (Section 1, Footer index 1 / wdHeaderFooterPrimary,
# stored as word/footer2.xml in the package) needs three edits:
#   1. Drop the stray "_GoBack" bookmark that sat in the
#      "...and Low Carbon Innovation" paragraph.
#   2. Rename the division in the "Mines ... Division" paragraph from
#      "Mines and Mineral Resources Division" to
#      "Mines, Competitiveness and Authorizations Division", and pick up
#      the paragraph-mark run formatting (Calibri instead of Arial) that
#      comes with it.
#   3. Remove the whole "Fax: {d.rc_office_fax_number}" paragraph and
#      move the "_GoBack" bookmark onto the paragraph that follows it
#      (the "Email: ..." paragraph).
#
# Range.Find/Replace and direct Range.Start/End slicing are not reliable
# against header/footer stories in this host, so the footer is rebuilt in
# one shot: clear the footer story, then re-insert the fully-formed target
# markup (paragraphs, table, bookmarks and all) via Range.InsertXML.

$targetFooterXml = @'
<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:wpc="http://schemas.microsoft.com/office/word/2010/wordprocessingCanvas" xmlns:cx="http://schemas.microsoft.com/office/drawing/2014/chartex" xmlns:cx1="http://schemas.microsoft.com/office/drawing/2015/9/8/chartex" xmlns:cx2="http://schemas.microsoft.com/office/drawing/2015/10/21/chartex" xmlns:cx3="http://schemas.microsoft.com/office/drawing/2016/5/9/chartex" xmlns:cx4="http://schemas.microsoft.com/office/drawing/2016/5/10/chartex" xmlns:cx5="http://schemas.microsoft.com/office/drawing/2016/5/11/chartex" xmlns:cx6="http://schemas.microsoft.com/office/drawing/2016/5/12/chartex" xmlns:cx7="http://schemas.microsoft.com/office/drawing/2016/5/13/chartex" xmlns:cx8="http://schemas.microsoft.com/office/drawing/2016/5/14/chartex" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:aink="http://schemas.microsoft.com/office/drawing/2016/ink" xmlns:am3d="http://schemas.microsoft.com/office/drawing/2017/model3d" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w15="http://schemas.microsoft.com/office/word/2012/wordml" xmlns:w16cex="http://schemas.microsoft.com/office/word/2018/wordml/cex" xmlns:w16cid="http://schemas.microsoft.com/office/word/2016/wordml/cid" xmlns:w16="http://schemas.microsoft.com/office/word/2018/wordml" xmlns:w16se="http://schemas.microsoft.com/office/word/2015/wordml/symex" xmlns:wpg="http://schemas.microsoft.com/office/word/2010/wordprocessingGroup" xmlns:wpi="http://schemas.microsoft.com/office/word/2010/wordprocessingInk" xmlns:wne="http://schemas.microsoft.com/office/word/2006/wordml" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" mc:Ignorable="w14 w15 w16se w16cid w16 w16cex wp14"><w:body><w:p w14:paraId="73F42B23" w14:textId="77777777" w:rsidR="00832E7B" w:rsidRPr="00832E7B" w:rsidRDefault="00832E7B" w:rsidP="00832E7B"><w:pPr><w:ind w:left="180" w:right="-288" w:hanging="270"/><w:textAlignment w:val="auto"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:u w:val="single"/><w:lang w:val="en-CA"/></w:rPr></w:pPr><w:r w:rsidRPr="00832E7B"><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:u w:val="single"/><w:lang w:val="en-CA"/></w:rPr><w:tab/></w:r><w:r w:rsidRPr="00832E7B"><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:u w:val="single"/><w:lang w:val="en-CA"/></w:rPr><w:tab/></w:r><w:r w:rsidRPr="00832E7B"><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:u w:val="single"/><w:lang w:val="en-CA"/></w:rPr><w:tab/></w:r><w:r w:rsidRPr="00832E7B"><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:u w:val="single"/><w:lang w:val="en-CA"/></w:rPr><w:tab/></w:r><w:r w:rsidRPr="00832E7B"><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:u w:val="single"/><w:lang w:val="en-CA"/></w:rPr><w:tab/></w:r><w:r w:rsidRPr="00832E7B"><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:u w:val="single"/><w:lang w:val="en-CA"/></w:rPr><w:tab/></w:r><w:r w:rsidRPr="00832E7B"><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:u w:val="single"/><w:lang w:val="en-CA"/></w:rPr><w:tab/></w:r><w:r w:rsidRPr="00832E7B"><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:u w:val="single"/><w:lang w:val="en-CA"/></w:rPr><w:tab/></w:r><w:r w:rsidRPr="00832E7B"><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:u w:val="single"/><w:lang w:val="en-CA"/></w:rPr><w:tab/></w:r><w:r w:rsidRPr="00832E7B"><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:u w:val="single"/><w:lang w:val="en-CA"/></w:rPr><w:tab/></w:r><w:r w:rsidRPr="00832E7B"><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:u w:val="single"/><w:lang w:val="en-CA"/></w:rPr><w:tab/></w:r><w:r w:rsidRPr="00832E7B"><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:u w:val="single"/><w:lang w:val="en-CA"/></w:rPr><w:tab/></w:r><w:r w:rsidRPr="00832E7B"><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:u w:val="single"/><w:lang w:val="en-CA"/></w:rPr><w:tab/></w:r><w:r w:rsidRPr="00832E7B"><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:u w:val="single"/><w:lang w:val="en-CA"/></w:rPr><w:tab/></w:r><w:r w:rsidRPr="00832E7B"><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:u w:val="single"/><w:lang w:val="en-CA"/></w:rPr><w:tab/></w:r></w:p><w:tbl><w:tblPr><w:tblW w:w="5000" w:type="pct"/><w:tblCellMar><w:left w:w="142" w:type="dxa"/><w:right w:w="142" w:type="dxa"/></w:tblCellMar><w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/></w:tblPr><w:tblGrid><w:gridCol w:w="3451"/><w:gridCol w:w="3563"/><w:gridCol w:w="3350"/></w:tblGrid><w:tr w:rsidR="00832E7B" w:rsidRPr="00832E7B" w14:paraId="6FBC54C2" w14:textId="77777777" w:rsidTr="00614FC2"><w:tc><w:tcPr><w:tcW w:w="1665" w:type="pct"/><w:shd w:val="clear" w:color="auto" w:fill="auto"/></w:tcPr><w:p w14:paraId="1CC35128" w14:textId="77777777" w:rsidR="000C64E1" w:rsidRDefault="00832E7B" w:rsidP="00832E7B"><w:pPr><w:tabs><w:tab w:val="center" w:pos="4320"/><w:tab w:val="right" w:pos="8640"/></w:tabs><w:ind w:left="293"/><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:b/><w:bCs/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr><w:bookmarkStart w:id="1" w:name="_Hlk37327873"/><w:r w:rsidRPr="00832E7B"><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:b/><w:bCs/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>Ministry of Energy, Mines</w:t></w:r></w:p><w:p w14:paraId="7927C9FC" w14:textId="22C5E3DD" w:rsidR="00832E7B" w:rsidRPr="00832E7B" w:rsidRDefault="000360D1" w:rsidP="00832E7B"><w:pPr><w:tabs><w:tab w:val="center" w:pos="4320"/><w:tab w:val="right" w:pos="8640"/></w:tabs><w:ind w:left="293"/><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:b/><w:bCs/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:b/><w:bCs/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>and Low Carbon Innovation</w:t></w:r></w:p><w:p w14:paraId="067CB486" w14:textId="77777777" w:rsidR="00832E7B" w:rsidRPr="00832E7B" w:rsidRDefault="00832E7B" w:rsidP="00832E7B"><w:pPr><w:tabs><w:tab w:val="left" w:pos="720"/><w:tab w:val="center" w:pos="4320"/><w:tab w:val="right" w:pos="8640"/></w:tabs><w:ind w:left="293"/><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:bCs/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:bCs/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>Mines, Competitiveness and Authorizations Division</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="0" w:type="auto"/><w:shd w:val="clear" w:color="auto" w:fill="auto"/></w:tcPr><w:p w14:paraId="79262269" w14:textId="77777777" w:rsidR="00832E7B" w:rsidRPr="00832E7B" w:rsidRDefault="00832E7B" w:rsidP="00832E7B"><w:pPr><w:tabs><w:tab w:val="left" w:pos="720"/><w:tab w:val="center" w:pos="4320"/><w:tab w:val="right" w:pos="8640"/></w:tabs><w:ind w:left="27" w:right="-548"/><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:b/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr><w:r w:rsidRPr="00832E7B"><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:b/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>Contact:</w:t></w:r></w:p><w:p w14:paraId="47031D04" w14:textId="77777777" w:rsidR="00832E7B" w:rsidRPr="00832E7B" w:rsidRDefault="00832E7B" w:rsidP="00832E7B"><w:pPr><w:tabs><w:tab w:val="left" w:pos="720"/><w:tab w:val="center" w:pos="4320"/><w:tab w:val="right" w:pos="8640"/></w:tabs><w:ind w:left="27" w:right="-548"/><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr><w:r w:rsidRPr="00832E7B"><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>Phone: {d.rc_office_phone_number}</w:t></w:r></w:p><w:p w14:paraId="1EA29E06" w14:textId="77777777" w:rsidR="00832E7B" w:rsidRPr="00832E7B" w:rsidRDefault="00832E7B" w:rsidP="00832E7B"><w:pPr><w:tabs><w:tab w:val="center" w:pos="4320"/><w:tab w:val="right" w:pos="8640"/></w:tabs><w:ind w:left="27" w:right="-169"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-CA"/></w:rPr></w:pPr><w:bookmarkStart w:id="2" w:name="_GoBack"/><w:bookmarkEnd w:id="2"/><w:r w:rsidRPr="00832E7B"><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>Email: {d.rc_office_email}</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="0" w:type="auto"/><w:shd w:val="clear" w:color="auto" w:fill="auto"/></w:tcPr><w:p w14:paraId="5F951562" w14:textId="77777777" w:rsidR="00832E7B" w:rsidRPr="00832E7B" w:rsidRDefault="00832E7B" w:rsidP="00832E7B"><w:pPr><w:tabs><w:tab w:val="left" w:pos="720"/><w:tab w:val="center" w:pos="4320"/><w:tab w:val="right" w:pos="8640"/></w:tabs><w:ind w:left="-533" w:right="-169" w:firstLine="425"/><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:b/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr><w:r w:rsidRPr="00832E7B"><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:b/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>Mailing Address:</w:t></w:r></w:p><w:p w14:paraId="75623C32" w14:textId="77777777" w:rsidR="00832E7B" w:rsidRPr="00832E7B" w:rsidRDefault="00832E7B" w:rsidP="00832E7B"><w:pPr><w:tabs><w:tab w:val="left" w:pos="876"/><w:tab w:val="center" w:pos="4320"/><w:tab w:val="right" w:pos="8640"/></w:tabs><w:ind w:left="-533" w:right="-169" w:firstLine="425"/><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr><w:r w:rsidRPr="00832E7B"><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>{d.rc_office_mailing_address_line_1}</w:t></w:r></w:p><w:p w14:paraId="54AC4412" w14:textId="77777777" w:rsidR="00832E7B" w:rsidRPr="00832E7B" w:rsidRDefault="00832E7B" w:rsidP="00832E7B"><w:pPr><w:tabs><w:tab w:val="left" w:pos="876"/><w:tab w:val="center" w:pos="4320"/><w:tab w:val="right" w:pos="8640"/></w:tabs><w:ind w:left="-533" w:right="-169" w:firstLine="425"/><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr><w:r w:rsidRPr="00832E7B"><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>{d.rc_office_mailing_address_line_2}</w:t></w:r></w:p></w:tc></w:tr><w:bookmarkEnd w:id="1"/></w:tbl><w:p w14:paraId="1831829E" w14:textId="77777777" w:rsidR="00832E7B" w:rsidRDefault="00832E7B"><w:pPr><w:pStyle w:val="Footer"/></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$d = $word.ActiveDocument
$ftr = $d.Sections(1).Footers(1)

# Clear the existing footer story content (collapses it to a single empty
# paragraph) so the subsequent InsertXML rebuilds the footer from scratch
# instead of being prepended ahead of the old content.
$rng = $ftr.Range
$rng.Text = ""

# Re-fetch the (now empty) footer range and inject the target OOXML.
$rng2 = $ftr.Range
$rng2.InsertXML($targetFooterXml)

$final = $ftr.Range
Write-Host "Updated footer text:" $final.Text
